$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dog")

# --- Row 20 ---
$ws.Range("A20").Value = 45810
$ws.Range("A20").NumberFormat = "m/d/yy"

$ws.Range("B20").Value = "PRESENCE"

$ws.Range("C20").Value = 0.39583333333333331
$ws.Range("C20").NumberFormat = "h:mm"
$ws.Range("D20").Value = 0.51388888888888884
$ws.Range("D20").NumberFormat = "h:mm"

$ws.Range("E20").Value = 13
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = "Sunny, mild"
$ws.Range("H20").Value = $true
$ws.Range("I20").Value = "3 minutes 50 seconds"
$ws.Range("J20").Value = 230
$ws.Range("K20").Value = "Primary sweeps"
$ws.Range("L20").Value = "Worked uphill, but the wind was more of a cross breeze."

# --- Row 21 ---
$ws.Range("A21").Value = 45810
$ws.Range("A21").NumberFormat = "m/d/yy"

$ws.Range("B21").Value = "PRESENCE"

$ws.Range("C21").Value = 0.53125
$ws.Range("C21").NumberFormat = "h:mm"
$ws.Range("D21").Value = 0.63541666666666663
$ws.Range("D21").NumberFormat = "h:mm"

$ws.Range("E21").Value = 11
$ws.Range("F21").Value = 9
$ws.Range("G21").Value = "Sunny, cool"
$ws.Range("H21").Value = $true
$ws.Range("I21").Value = "7 minutes 41 seconds"
$ws.Range("J21").Value = 461
$ws.Range("K21").Value = "Secondary checks"
$ws.Range("L21").Value = "Worked downhill. Missed on transects (apparently because I was standing on it when we passed) but picked up during the perimeter sweep."

# --- View changes ---
$ws.Range("L28").Select()
